$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 6 new localization rows (56-61) -----------------------------

# Row 56: lang_next
$ws.Range("A56").Value = 'lang_next'
$ws.Range("B56").Value = 'Tiếp Theo'
$ws.Range("C56").Value = 'Next'

# Row 57: lang_partner_list
$ws.Range("A57").Value = 'lang_partner_list'
$ws.Range("B57").Value = 'Người Đi Cùng'
$ws.Range("C57").Value = 'Partner'

# Row 58: lang_pick_method_with_parent
$ws.Range("A58").Value = 'lang_pick_method_with_parent'
$ws.Range("B58").Value = 'PH giao HS tận tay Giám sát Xe'
$ws.Range("C58").Value = 'Parent go with Student until kick off'

# Row 59: lang_pick_method_by_student
$ws.Range("A59").Value = 'lang_pick_method_by_student'
$ws.Range("B59").Value = 'HS tự đón xe và về tại điểm đón trả'
$ws.Range("C59").Value = 'Student kickk off and go off istself'

# Row 60: lang_service_start_date
$ws.Range("A60").Value = 'lang_service_start_date'
$ws.Range("B60").Value = 'Ngày bắt đầu dịch vụ'
$ws.Range("C60").Value = 'Service Starting Date '

# Row 61: lang_alert_wrong_year
$ws.Range("A61").Value = 'lang_alert_wrong_year'
$ws.Range("B61").Value = 'Vui lòng điều chỉnh năm đăng kí @year@'
$ws.Range("C61").Value = 'Please adjust Year @year@'

# --- Match formatting of the surrounding rows ----------------------------
# Column A uses the same "Consolas" style as row 55 for rows 56-59 ...
$ws.Range("A55").Copy()
$ws.Range("A56:A59").PasteSpecial(-4122)

# ... and the same "Consolas/orange" style as row 54 for rows 60-61.
$ws.Range("A54").Copy()
$ws.Range("A60:A61").PasteSpecial(-4122)

# Columns B and C keep the plain "Calibri" style used throughout column B/C.
$ws.Range("B55").Copy()
$ws.Range("B56:C61").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Sheet-level cosmetics to mirror the authored workbook ---------------
$ws.Columns.Item(1).ColumnWidth = 50.25

$ws.Range("C61").Select() | Out-Null
